$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 416.33334
$ws.Range("I2").Value = 380.875
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 380.875
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -267.875
$ws.Range("N2").Value = -926
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("H40").Value = 2450
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H88").Value = 22730136
$ws.Range("J88").Value = 2933.1333
$ws.Range("L88").Value = 2933.1333
$ws.Range("N88").Value = -3745.1333
$ws.Range("H91").Value = 22730136
$ws.Range("J91").Value = 2933.1333
$ws.Range("L91").Value = 2933.1333
$ws.Range("N91").Value = -5741.1333
$ws.Range("H92").Value = 20834112
$ws.Range("I92").Value = 27778274
$ws.Range("K92").Value = 27778274
$ws.Range("M92").Value = -27777026
$ws.Range("H94").Value = 2379.4
$ws.Range("I94").Value = 2474.25
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 2474.25
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -2023.25
$ws.Range("N94").Value = -2902
$ws.Range("H99").Value = 1251152
$ws.Range("I99").Value = 551
$ws.Range("K99").Value = 1653
$ws.Range("M99").Value = -155
$ws.Range("H137").Value = 38404.223
$ws.Range("I137").Value = 1429.6923
$ws.Range("J137").Value = 72737.71000000001
$ws.Range("K137").Value = 4289.0769
$ws.Range("L137").Value = 218213.13
$ws.Range("M137").Value = -1739.0769
$ws.Range("N137").Value = -223313.13
$ws.Range("H138").Value = 5041.2856
$ws.Range("J138").Value = 4730.3
$ws.Range("L138").Value = 14190.9
$ws.Range("N138").Value = -24470.9
$ws.Range("H141").Value = 2423.3076
$ws.Range("I141").Value = 2041.9166
$ws.Range("K141").Value = 6125.7498
$ws.Range("M141").Value = -945.7497999999996
$ws.Range("M5").ClearContents()
$ws.Range("N40").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 18850.902
$ws.Range("I61").Value = 23847.527
$ws.Range("K61").Value = 23847.527
$ws.Range("M61").Value = -23635.527
$ws.Range("H74").Value = 1916
$ws.Range("I74").Value = 678
$ws.Range("K74").Value = 678
$ws.Range("M74").Value = 196
$ws.Range("H77").Value = 1916
$ws.Range("I77").Value = 678
$ws.Range("K77").Value = 3390
$ws.Range("M77").Value = 978
$ws.Range("H136").Value = 18850.902
$ws.Range("I136").Value = 23847.527
$ws.Range("K136").Value = 71542.58099999999
$ws.Range("M136").Value = -68992.58099999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8461.727999999999
$ws.Range("I134").Value = 8845.157999999999
$ws.Range("K134").Value = 26535.474
$ws.Range("M134").Value = -24000.474

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1183.4286
$ws.Range("I5").Value = 1294.6666
$ws.Range("K5").Value = 1294.6666
$ws.Range("M5").Value = -1182.6666
$ws.Range("H6").Value = 9969
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("H7").Value = 24.8125
$ws.Range("I7").Value = 22.466667
$ws.Range("K7").Value = 22.466667
$ws.Range("M7").Value = 90.533333
$ws.Range("H31").Value = 2537.8823
$ws.Range("I31").Value = 1796.2142
$ws.Range("K31").Value = 1796.2142
$ws.Range("M31").Value = -1501.2142
$ws.Range("H34").Value = 2537.8823
$ws.Range("I34").Value = 1796.2142
$ws.Range("K34").Value = 1796.2142
$ws.Range("M34").Value = -1594.2142
$ws.Range("H50").Value = 13280
$ws.Range("J50").Value = 13280
$ws.Range("L50").Value = 13280
$ws.Range("N50").Value = -14530
$ws.Range("H59").Value = 30690.8
$ws.Range("J59").Value = 17800
$ws.Range("L59").Value = 17800
$ws.Range("N59").Value = -20090
$ws.Range("H60").Value = 26220.77
$ws.Range("J60").Value = 26220.77
$ws.Range("L60").Value = 26220.77
$ws.Range("N60").Value = -27242.77
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N141").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 129.04546
$ws.Range("J12").Value = 137.06667
$ws.Range("L12").Value = 411.20001
$ws.Range("N12").Value = -757.20001
$ws.Range("H122").Value = 1949.5
$ws.Range("J122").Value = 1949.5
$ws.Range("L122").Value = 17545.5
$ws.Range("N122").Value = -22445.5
$ws.Range("H129").Value = 61590.582
$ws.Range("I129").Value = 699.6667
$ws.Range("J129").Value = 81887.55499999999
$ws.Range("K129").Value = 2099.0001
$ws.Range("L129").Value = 245662.665
$ws.Range("M129").Value = 2900.9999
$ws.Range("N129").Value = -255662.665
$ws.Range("H131").Value = 12857.061
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 12857.061
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 38571.183
$ws.Range("N131").Value = -48651.183
$ws.Range("H137").Value = 7206.375
$ws.Range("I137").Value = 3523.1667
$ws.Range("K137").Value = 10569.5001
$ws.Range("M137").Value = -5469.500100000001
$ws.Range("M131").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2595.2
$ws.Range("I80").Value = 2749
$ws.Range("K80").Value = 2749
$ws.Range("M80").Value = -1751
$ws.Range("H83").Value = 2595.2
$ws.Range("I83").Value = 2749
$ws.Range("K83").Value = 13745
$ws.Range("M83").Value = -8753

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3559.8
$ws.Range("J46").Value = 3724.75
$ws.Range("L46").Value = 3724.75
$ws.Range("N46").Value = -4100.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 15878.6
$ws.Range("I96").Value = 1797.6666
$ws.Range("K96").Value = 1797.6666
$ws.Range("M96").Value = -424.6666
$ws.Range("H122").Value = 76869.09
$ws.Range("I122").Value = 116308.86
$ws.Range("J122").Value = 7849.5
$ws.Range("K122").Value = 348926.58
$ws.Range("L122").Value = 23548.5
$ws.Range("M122").Value = -346476.58
$ws.Range("N122").Value = -28448.5
$ws.Range("H126").Value = 4393.4546
$ws.Range("I126").Value = 5595.3335
$ws.Range("J126").Value = 2951.2
$ws.Range("K126").Value = 16786.0005
$ws.Range("L126").Value = 8853.599999999999
$ws.Range("M126").Value = -14316.0005
$ws.Range("N126").Value = -13793.6
$ws.Range("H132").Value = 1359.2
$ws.Range("I132").Value = 1010.58826
$ws.Range("K132").Value = 3031.76478
$ws.Range("M132").Value = -501.76478
